$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "added objects for trials" - a new (currently blank) leading column is
# introduced before the existing "Case ID" column, pushing the whole
# header row right by one (Case ID -> B, Trial Code -> C, ... Ethnicity -> I)
# and leaving a new trailing "Ethnicity"-style object slot at I1.
$ws.Columns("A:A").Insert()

# The single data/record row loses its populated trial values (B2:I2) -
# only the leading placeholder column (A2) remains.
$ws.Range("B2:I2").ClearContents()

# Touch the new leading column's cells with a no-op format write so the
# engine keeps them as present-but-empty cells (matching row/column shape)
# instead of dropping them from the sparse sheet entirely.
$ws.Range("A1:A2").Font.Bold = $false
